$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "<wat>"
$ws.Range("C2").Value = 50

$ws.Range("B3").Value = "<tab>"
$ws.Range("C3").Value = 51

$ws.Range("B4").Value = "<but>"
$ws.Range("C4").Value = 50

$ws.Range("C5").Value = 49

$ws.Range("B6").Value = "<these>"
$ws.Range("C6").Value = 52

$ws.Range("B7").Value = "<on>"

$ws.Range("C8").Value = 47

$ws.Range("B9").Value = "<delete>"
$ws.Range("C9").Value = 46

$ws.Range("C10").Value = 48

$ws.Range("C11").Value = 46

$ws.Range("B12").Value = "<down>"
$ws.Range("C12").Value = 46

$ws.Range("B13").Value = "<for>"
$ws.Range("C13").Value = 43

$ws.Range("B14").Value = "<by>"
$ws.Range("C14").Value = 47

$ws.Range("C15").Value = 49

$ws.Range("C16").Value = 46

$ws.Range("C18").Value = 48
